$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 354 (the "カルマン線" post), shifting subsequent rows up by one.
$ws.Rows.Item(354).Delete()
